# Update ADANI closing-date values that were wrong for the algo.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADANI")

# Row 7 - summary row (Open/High/Low/Close/Prev)
$ws.Range("F7").Value = 2801.7
$ws.Range("G7").Value = 2821.95
$ws.Range("H7").Value = 2767.1
$ws.Range("I7").Value = 2799.75
$ws.Range("J7").Value = 2783.85

# Row 9
$ws.Range("G9").Value = 2817.8
$ws.Range("H9").Value = 2725
$ws.Range("I9").Value = 2782.9

# Row 10
$ws.Range("G10").Value = 2810
$ws.Range("H10").Value = 2772
$ws.Range("I10").Value = 2802

# Row 11
$ws.Range("G11").Value = 2817.9
$ws.Range("H11").Value = 2779.75
$ws.Range("I11").Value = 2782.85

# Row 12
$ws.Range("G12").Value = 2798.15
$ws.Range("H12").Value = 2767.1
$ws.Range("I12").Value = 2795.4

# Row 13
$ws.Range("G13").Value = 2798.8
$ws.Range("H13").Value = 2780
$ws.Range("I13").Value = 2787.3

# Row 14
$ws.Range("G14").Value = 2799
$ws.Range("H14").Value = 2786.15
$ws.Range("I14").Value = 2794.5

# Row 15
$ws.Range("G15").Value = 2805.95
$ws.Range("H15").Value = 2789.05
$ws.Range("I15").Value = 2795.4

# Row 16
$ws.Range("G16").Value = 2800.8
$ws.Range("H16").Value = 2789
$ws.Range("I16").Value = 2791.9

# Row 17
$ws.Range("G17").Value = 2798.9
$ws.Range("H17").Value = 2785.1
$ws.Range("I17").Value = 2798.4

# Row 18
$ws.Range("G18").Value = 2803.95
$ws.Range("H18").Value = 2773
$ws.Range("I18").Value = 2784.4

# Row 19
$ws.Range("G19").Value = 2791.95
$ws.Range("H19").Value = 2777.1
$ws.Range("I19").Value = 2778.95

# Row 20
$ws.Range("G20").Value = 2787.45
$ws.Range("H20").Value = 2774.2
$ws.Range("I20").Value = 2776.85

# Row 21
$ws.Range("G21").Value = 2821.95
$ws.Range("H21").Value = 2774
$ws.Range("I21").Value = 2814
